$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.185.78'
$ws.Range('E2').Value = '  -3.75%  '
$ws.Range('D3').Value = '3.646.70'
$ws.Range('E3').Value = '  -5.29%  '
$ws.Range('E4').Value = '  +0.23%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '587.98'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.69%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '178.62'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +5.92%  '
$ws.Range('D7').Value = '3.640.24'
$ws.Range('E7').Value = '  -5.24%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.626'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -5.72%  '
$ws.Range('E9').Value = '  +0.07%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.707'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -5.08%  '
$ws.Range('E11').Value = '  -8.91%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '55.72'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +5.12%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000289'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -9.61%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '10.52'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -7.09%  '
$ws.Range('D15').Value = '4.246.39'
$ws.Range('E15').Value = '  -4.56%  '
$ws.Range('D16').Value = '3.658.85'
$ws.Range('E16').Value = '  -5.03%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '19.18'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -8.32%  '
$ws.Range('E18').Value = '  -2.40%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.69'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -8.34%  '
$ws.Range('E20').Value = '  -7.21%  '
$ws.Range('D21').Value = '68.020.31'
$ws.Range('E21').Value = '  -3.64%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '407.14'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -6.50%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.52'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -5.62%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '88.10'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -5.81%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.00'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -8.70%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '4.07'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '12.63'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -8.60%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '10.72'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -8.56%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.04'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.32%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.40'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -10.20%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '32.45'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -7.23%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '7.15'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -14.08%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '12.27'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -8.86%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.116'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -6.97%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '64.45'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -6.15%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '42.83'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -11.06%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '589.55'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -8.53%  '
$ws.Range('D38').Value = '0.0₃0879'
$ws.Range('E38').Value = '  -11.27%  '
$ws.Range('E39').Value = '  -0.31%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.395'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -8.56%  '
$ws.Range('E41').Value = '  +0.13%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.136'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -7.50%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.00'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -7.05%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.68'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -9.23%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0434'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -7.25%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.83'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -12.91%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.134'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -6.50%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.69'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -3.22%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '8.95'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -10.20%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '3.14'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('D51').Value = '2.692.11'
$ws.Range('E51').Value = '  -7.58%  '
